# Updates cryptos list values (price + 1h volume change) per commit
# "Updated cryptos list on Wed Aug 28 17:10:16 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.952.09'
$ws.Range("E2").Value = '  -4.31%  '
$ws.Range("D3").Value = '''2.484.53'
$ws.Range("E3").Value = '  -3.29%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '''534.75'
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").Value = '''142.45'
$ws.Range("E6").Value = '  -7.51%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''0.571'
$ws.Range("E8").Value = '  -4.50%  '
$ws.Range("D9").Value = '''2.510.64'
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("D10").Value = '''0.0997'
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("D12").Value = '''5.50'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").Value = '''0.350'
$ws.Range("E13").Value = '  -3.48%  '
$ws.Range("D14").Value = '''2.952.42'
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").Value = '''23.78'
$ws.Range("E15").Value = '  -6.33%  '
$ws.Range("D16").Value = '''58.922.82'
$ws.Range("E16").Value = '  -4.19%  '
$ws.Range("D17").Value = '''0.0000138'
$ws.Range("E17").Value = '  -3.81%  '
$ws.Range("D18").Value = '''2.494.27'
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("D19").Value = '''11.33'
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").Value = '''4.28'
$ws.Range("E20").Value = '  -5.29%  '
$ws.Range("D21").Value = '''322.83'
$ws.Range("E21").Value = '  -3.93%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '''5.74'
$ws.Range("E23").Value = '  -4.90%  '
$ws.Range("D24").Value = '''60.77'
$ws.Range("E24").Value = '  -3.72%  '
$ws.Range("D25").Value = '''0.436'
$ws.Range("E25").Value = '  -11.56%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").Value = '''0.161'
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("D28").Value = '''2.596.83'
$ws.Range("E28").Value = '  -3.44%  '
$ws.Range("D29").Value = '''7.74'
$ws.Range("E29").Value = '  -4.05%  '
$ws.Range("D30").Value = '''6.86'
$ws.Range("E30").Value = '  -7.33%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '''1.25'
$ws.Range("E31").Value = '  -2.17%  '
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '''0.0₃0764'
$ws.Range("E32").Value = '  -8.30%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.77'
$ws.Range("E33").Value = '  -6.23%  '
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").Value = '''156.91'
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").Value = '''1.41'
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").Value = '''18.53'
$ws.Range("E37").Value = '  -3.11%  '
$ws.Range("D38").Value = '''4.37'
$ws.Range("E38").Value = '  -7.12%  '
$ws.Range("D39").Value = '''1.60'
$ws.Range("E39").Value = '  -10.14%  '
$ws.Range("D40").Value = '''5.87'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = '''308.88'
$ws.Range("E41").Value = '  -7.45%  '
$ws.Range("D42").Value = '''36.67'
$ws.Range("E42").Value = '  -2.02%  '
$ws.Range("D43").Value = '''3.65'
$ws.Range("E43").Value = '  -7.48%  '
$ws.Range("D44").Value = '''0.793'
$ws.Range("E44").Value = '  -15.05%  '
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("D46").Value = '''0.595'
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("D48").Value = '''124.66'
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("D49").Value = '''0.0925'
$ws.Range("E49").Value = '  -4.02%  '
$ws.Range("D50").Value = '''18.56'
$ws.Range("E50").Value = '  -4.54%  '
$ws.Range("D51").Value = '''0.0517'
$ws.Range("E51").Value = '  -4.88%  '
